$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "Op-Amps" worksheet before "LEDs"
# ---------------------------------------------------------------
$leds = $wb.Worksheets.Item("LEDs")
$opamps = $wb.Worksheets.Add($leds)
$opamps.Name = "Op-Amps"

$headers = @("id","Symbol","Footprint","MPN","Function","Channel Count","Gain","Bandwidth","Slew Rate","Input Offset Voltage","Input Bias Current","Pin Count","Manufacturer","Datasheet","Description","Distributer PN","Price","Footprint Filters","Keywords","No BOM","Schematic Only")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $opamps.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$opamps.Cells.Item(2, 1).Value = "primary"
$opamps.Cells.Item(3, 1).Value = "auto"

$opamps.Range("A1:U3").Columns.AutoFit() | Out-Null

# ---------------------------------------------------------------
# 2. Update view/selection state on a couple of existing sheets
# ---------------------------------------------------------------
$capacitors = $wb.Worksheets.Item("Capacitors")
$capacitors.Range("X23").Select()

$diodes = $wb.Worksheets.Item("Diodes")
$diodes.Range("R16").Select()

# Make Op-Amps the active sheet / selected cell last, so it ends up
# as the active tab in the saved workbook.
$opamps.Activate()
$opamps.Range("A4").Select()
